$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.185.96"
$ws.Cells.Item(2, 5).Value = "  +4.22%  "

$ws.Cells.Item(3, 4).Value = "2.950.00"
$ws.Cells.Item(3, 5).Value = "  +1.97%  "

$ws.Cells.Item(4, 5).Value = "  -0.20%  "

$ws.Cells.Item(5, 4).Value = "577.16"
$ws.Cells.Item(5, 5).Value = "  +0.82%  "

$ws.Cells.Item(6, 4).Value = "151.51"
$ws.Cells.Item(6, 5).Value = "  +5.09%  "

$ws.Cells.Item(7, 5).Value = "  +0.07%  "

$ws.Cells.Item(8, 4).Value = "2.947.29"
$ws.Cells.Item(8, 5).Value = "  +1.78%  "

$ws.Cells.Item(9, 4).Value = "0.508"
$ws.Cells.Item(9, 5).Value = "  +0.90%  "

$ws.Cells.Item(10, 5).Value = "  +4.52%  "

$ws.Cells.Item(11, 4).Value = "0.151"
$ws.Cells.Item(11, 5).Value = "  +1.97%  "

$ws.Cells.Item(12, 4).Value = "0.443"
$ws.Cells.Item(12, 5).Value = "  +2.68%  "

$ws.Cells.Item(13, 4).Value = "0.0000240"
$ws.Cells.Item(13, 5).Value = "  +3.07%  "

$ws.Cells.Item(14, 4).Value = "34.09"
$ws.Cells.Item(14, 5).Value = "  +5.71%  "

$ws.Cells.Item(15, 5).Value = "  +0.70%  "

$ws.Cells.Item(16, 2).Value = "WrappedBTC"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(16, 4).Value = "64.083.20"
$ws.Cells.Item(16, 5).Value = "  +3.92%  "

$ws.Cells.Item(17, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(17, 4).Value = "3.440.84"
$ws.Cells.Item(17, 5).Value = "  +2.01%  "

$ws.Cells.Item(18, 4).Value = "6.86"
$ws.Cells.Item(18, 5).Value = "  +3.87%  "

$ws.Cells.Item(19, 4).Value = "2.950.71"
$ws.Cells.Item(19, 5).Value = "  +2.10%  "

$ws.Cells.Item(20, 4).Value = "444.30"
$ws.Cells.Item(20, 5).Value = "  +2.68%  "

$ws.Cells.Item(21, 4).Value = "13.42"
$ws.Cells.Item(21, 5).Value = "  +1.30%  "

$ws.Cells.Item(22, 4).Value = "0.671"
$ws.Cells.Item(22, 5).Value = "  +2.47%  "

$ws.Cells.Item(23, 4).Value = "7.11"
$ws.Cells.Item(23, 5).Value = "  +3.01%  "

$ws.Cells.Item(24, 4).Value = "80.08"
$ws.Cells.Item(24, 5).Value = "  +0.68%  "

$ws.Cells.Item(25, 4).Value = "10.79"
$ws.Cells.Item(25, 5).Value = "  +6.94%  "

$ws.Cells.Item(26, 4).Value = "12.16"
$ws.Cells.Item(26, 5).Value = "  +2.89%  "

$ws.Cells.Item(27, 4).Value = "2.16"
$ws.Cells.Item(27, 5).Value = "  +6.31%  "

$ws.Cells.Item(28, 5).Value = "  -0.03%  "

$ws.Cells.Item(29, 2).Value = "PEPE"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(29, 4).Value = "0.0000110"
$ws.Cells.Item(29, 5).Value = "  +3.30%  "

$ws.Cells.Item(30, 2).Value = "NEARProtocol"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(30, 4).Value = "7.56"
$ws.Cells.Item(30, 5).Value = "  +7.98%  "

$ws.Cells.Item(31, 4).Value = "2.13"
$ws.Cells.Item(31, 5).Value = "  +2.86%  "

$ws.Cells.Item(32, 4).Value = "2.54"
$ws.Cells.Item(32, 5).Value = "  +0.86%  "

$ws.Cells.Item(33, 4).Value = "0.109"
$ws.Cells.Item(33, 5).Value = "  +2.33%  "

$ws.Cells.Item(34, 4).Value = "26.34"
$ws.Cells.Item(34, 5).Value = "  +3.14%  "

$ws.Cells.Item(35, 5).Value = "  -0.31%  "

$ws.Cells.Item(36, 4).Value = "0.967"
$ws.Cells.Item(36, 5).Value = "  +0.87%  "

$ws.Cells.Item(37, 2).Value = "Filecoin"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(37, 4).Value = "5.58"
$ws.Cells.Item(37, 5).Value = "  +2.70%  "

$ws.Cells.Item(38, 2).Value = "Stacks"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(38, 4).Value = "2.11"
$ws.Cells.Item(38, 5).Value = "  +8.22%  "

$ws.Cells.Item(39, 4).Value = "2.99"
$ws.Cells.Item(39, 5).Value = "  +0.09%  "

$ws.Cells.Item(40, 4).Value = "49.01"
$ws.Cells.Item(40, 5).Value = "  -0.10%  "

$ws.Cells.Item(41, 4).Value = "43.07"
$ws.Cells.Item(41, 5).Value = "  +12.29%  "

$ws.Cells.Item(42, 4).Value = "0.117"
$ws.Cells.Item(42, 5).Value = "  +1.90%  "

$ws.Cells.Item(43, 4).Value = "0.290"
$ws.Cells.Item(43, 5).Value = "  +8.19%  "

$ws.Cells.Item(44, 4).Value = "8.24"
$ws.Cells.Item(44, 5).Value = "  -0.19%  "

$ws.Cells.Item(45, 4).Value = "377.51"
$ws.Cells.Item(45, 5).Value = "  +10.78%  "

$ws.Cells.Item(46, 2).Value = "Maker"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(46, 4).Value = "2.739.64"
$ws.Cells.Item(46, 5).Value = "  +2.34%  "

$ws.Cells.Item(47, 2).Value = "VeChain"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(47, 4).Value = "0.0347"
$ws.Cells.Item(47, 5).Value = "  +4.27%  "

$ws.Cells.Item(48, 4).Value = "134.41"
$ws.Cells.Item(48, 5).Value = "  +0.77%  "

$ws.Cells.Item(49, 5).Value = "  +0.02%  "

$ws.Cells.Item(50, 2).Value = "FLOKI"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(50, 4).Value = "0.000218"
$ws.Cells.Item(50, 5).Value = "  +10.01%  "

$ws.Cells.Item(51, 2).Value = "Stellar"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(51, 4).Value = "0.105"
$ws.Cells.Item(51, 5).Value = "  +2.24%  "
